# Regenerate orders with updated distance/size codes.
#   Distance codes:  D80 -> D86,  D64 -> D69,  D51 -> D55
#   Size codes:       S30 -> S31  (S20/S25 unchanged)
#
# These substitutions touch every string cell that embeds a distance or
# size token: the Condition / Filename_Left / Filename_Right / Distance /
# Size columns, including the lookup rows further down the sheet (e.g.
# "D64", "D80", "D51", "S30" on their own).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value()
        if ($v -is [string]) {
            $nv = $v
            $nv = $nv.Replace("D80", "D86")
            $nv = $nv.Replace("D64", "D69")
            $nv = $nv.Replace("D51", "D55")
            $nv = $nv.Replace("S30", "S31")
            if ($nv -ne $v) {
                $cell.Value = $nv
            }
        }
    }
}
